$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.98080691254106
$ws.Range("C2").Value = 9.227977367057209
$ws.Range("E2").Value = 15.47374107811253
$ws.Range("F2").Value = 40.59432522638034
$ws.Range("G2").Value = 3.68166311293257
$ws.Range("I2").Value = 27.68757089880922
$ws.Range("J2").Value = 8.703359469589422
$ws.Range("K2").Value = 8.861881177781523
$ws.Range("L2").Value = 12.36240689968064
$ws.Range("M2").Value = 15.43118098463983
$ws.Range("O2").Value = 27.44555715225936

$ws.Range("B3").Value = 11.75606587258971
$ws.Range("C3").Value = 9.219855200171917
$ws.Range("E3").Value = 15.48799953983259
$ws.Range("F3").Value = 40.66007242091329
$ws.Range("G3").Value = 3.683339155719638
$ws.Range("I3").Value = 27.79224900545679
$ws.Range("J3").Value = 8.694431752776186
$ws.Range("K3").Value = 8.701703670363344
$ws.Range("L3").Value = 12.35868463020991
$ws.Range("M3").Value = 15.38844718724152
$ws.Range("O3").Value = 27.54547365894658

$ws.Range("B4").Value = 11.61758283061694
$ws.Range("C4").Value = 9.214958096058096
$ws.Range("E4").Value = 15.49866190033494
$ws.Range("F4").Value = 40.70872149474445
$ws.Range("G4").Value = 3.684423010456533
$ws.Range("I4").Value = 27.86090826929865
$ws.Range("J4").Value = 8.688922276077211
$ws.Range("K4").Value = 8.60293310515886
$ws.Range("L4").Value = 12.35789956903391
$ws.Range("M4").Value = 15.36403832197225
$ws.Range("O4").Value = 27.6116039807283

$ws.Range("B5").Value = 11.56109986107825
$ws.Range("C5").Value = 9.212985523928769
$ws.Range("E5").Value = 15.50348739415921
$ws.Range("F5").Value = 40.73062679535121
$ws.Range("G5").Value = 3.684878501336054
$ws.Range("I5").Value = 27.88999096175391
$ws.Range("J5").Value = 8.686670386181538
$ws.Range("K5").Value = 8.56262754429917
$ws.Range("L5").Value = 12.35795846500632
$ws.Range("M5").Value = 15.35455819938709
$ws.Range("O5").Value = 27.63975440118954

$ws.Range("B6").Value = 11.55172006441608
$ws.Range("C6").Value = 9.21265938351787
$ws.Range("E6").Value = 15.50431770682258
$ws.Range("F6").Value = 40.73438974625205
$ws.Range("G6").Value = 3.684954970653349
$ws.Range("I6").Value = 27.89488678274875
$ws.Range("J6").Value = 8.686296062519611
$ws.Range("K6").Value = 8.55593295254868
$ws.Range("L6").Value = 12.35799116616831
$ws.Range("M6").Value = 15.35301240356253
$ws.Range("O6").Value = 27.64450132735469

$ws.Range("B7").Value = 11.61682118577352
$ws.Range("C7").Value = 9.214931399423728
$ws.Range("E7").Value = 15.49872503209807
$ws.Range("F7").Value = 40.70900849639759
$ws.Range("G7").Value = 3.684429097388704
$ws.Range("I7").Value = 27.86129602005725
$ws.Range("J7").Value = 8.68889193322395
$ws.Range("K7").Value = 8.602389689096976
$ws.Range("L7").Value = 12.35789882759741
$ws.Range("M7").Value = 15.36390857172809
$ws.Range("O7").Value = 27.6119787618418

$ws.Range("B8").Value = 11.90346171635192
$ws.Range("C8").Value = 9.225158428685667
$ws.Range("E8").Value = 15.47826197792809
$ws.Range("F8").Value = 40.6152754812171
$ws.Range("G8").Value = 3.682229672859262
$ws.Range("I8").Value = 27.72275377556072
$ws.Range("J8").Value = 8.700287123995148
$ws.Range("K8").Value = 8.806769054852513
$ws.Range("L8").Value = 12.36081299238887
$ws.Range("M8").Value = 15.41607061906248
$ws.Range("O8").Value = 27.47901573688895

$ws.Range("B9").Value = 12.45853030850055
$ws.Range("C9").Value = 9.245909514757033
$ws.Range("E9").Value = 15.45323265986855
$ws.Range("F9").Value = 40.49722010019068
$ws.Range("G9").Value = 3.678349166561505
$ws.Range("I9").Value = 27.4858562987605
$ws.Range("J9").Value = 8.722401449595079
$ws.Range("K9").Value = 9.202076657309481
$ws.Range("L9").Value = 12.37836246859888
$ws.Range("M9").Value = 15.5325777509941
$ws.Range("O9").Value = 27.25623384211229

$ws.Range("B10").Value = 12.85801798860294
$ws.Range("C10").Value = 9.26155698284443
$ws.Range("E10").Value = 15.44399547495156
$ws.Range("F10").Value = 40.4506203048925
$ws.Range("G10").Value = 3.675759168229381
$ws.Range("I10").Value = 27.33298287295411
$ws.Range("J10").Value = 8.738494102822544
$ws.Range("K10").Value = 9.486410889079945
$ws.Range("L10").Value = 12.39837105168649
$ws.Range("M10").Value = 15.62641959739183
$ws.Range("O10").Value = 27.11572556669046

$ws.Range("B11").Value = 13.03715091912537
$ws.Range("C11").Value = 9.268757821876552
$ws.Range("E11").Value = 15.4417677511598
$ws.Range("F11").Value = 40.43813479186849
$ws.Range("G11").Value = 3.674637014800935
$ws.Range("I11").Value = 27.26802907801231
$ws.Range("J11").Value = 8.745778569144765
$ws.Range("K11").Value = 9.613890770939744
$ws.Range("L11").Value = 12.40899464744948
$ws.Range("M11").Value = 15.6708039690119
$ws.Range("O11").Value = 27.05684232049775

$ws.Range("B12").Value = 13.10454851136874
$ws.Range("C12").Value = 9.271496044719123
$ws.Range("E12").Value = 15.44120684028235
$ws.Range("F12").Value = 40.43465852841667
$ws.Range("G12").Value = 3.674220101204906
$ws.Range("I12").Value = 27.24409240147603
$ws.Range("J12").Value = 8.748531525660347
$ws.Range("K12").Value = 9.661853145740542
$ws.Range("L12").Value = 12.4132340703888
$ws.Range("M12").Value = 15.68784623593355
$ws.Range("O12").Value = 27.03526940338324

$ws.Range("B13").Value = 13.09005359506535
$ws.Range("C13").Value = 9.270905819070917
$ws.Range("E13").Value = 15.44131509064209
$ws.Range("F13").Value = 40.43535155872585
$ws.Range("G13").Value = 3.674309534853702
$ws.Range("I13").Value = 27.24921824261267
$ws.Range("J13").Value = 8.747938877422857
$ws.Range("K13").Value = 9.65153810768007
$ws.Range("L13").Value = 12.41231144409351
$ws.Range("M13").Value = 15.68416557773214
$ws.Range("O13").Value = 27.03988326284599

$ws.Range("B14").Value = 13.04270484822246
$ws.Range("C14").Value = 9.268982870949706
$ws.Range("E14").Value = 15.441715946683
$ws.Range("F14").Value = 40.43782372015402
$ws.Range("G14").Value = 3.674602554523439
$ws.Range("I14").Value = 27.26604656735886
$ws.Range("J14").Value = 8.746005168824452
$ws.Range("K14").Value = 9.617843146586589
$ws.Range("L14").Value = 12.40933910205222
$ws.Range("M14").Value = 15.67220139144886
$ws.Range("O14").Value = 27.05505297051434

$ws.Range("B15").Value = 13.01364376812258
$ws.Range("C15").Value = 9.26780648294411
$ws.Range("E15").Value = 15.44199825801
$ws.Range("F15").Value = 40.43950095629265
$ws.Range("G15").Value = 3.674783081042075
$ws.Range("I15").Value = 27.27644034354866
$ws.Range("J15").Value = 8.744819989194726
$ws.Range("K15").Value = 9.597162200686547
$ws.Range("L15").Value = 12.40754658050292
$ws.Range("M15").Value = 15.66490329965493
$ws.Range("O15").Value = 27.06443928367484

$ws.Range("B16").Value = 12.84625353744709
$ws.Range("C16").Value = 9.261088015278146
$ws.Range("E16").Value = 15.44418068754511
$ws.Range("F16").Value = 40.45161150321637
$ws.Range("G16").Value = 3.67583362817593
$ws.Range("I16").Value = 27.33732010061118
$ws.Range("J16").Value = 8.73801729150323
$ws.Range("K16").Value = 9.478038493549693
$ws.Range("L16").Value = 12.39770719894061
$ws.Range("M16").Value = 15.62355233217731
$ws.Range("O16").Value = 27.11967511353596

$ws.Range("B17").Value = 12.74285430426811
$ws.Range("C17").Value = 9.256987304966955
$ws.Range("E17").Value = 15.44602436356491
$ws.Range("F17").Value = 40.46127189364064
$ws.Range("G17").Value = 3.676492433130278
$ws.Range("I17").Value = 27.37584323339754
$ws.Range("J17").Value = 8.733834573610091
$ws.Range("K17").Value = 9.40445070504199
$ws.Range("L17").Value = 12.39205917425797
$ws.Range("M17").Value = 15.59861280027758
$ws.Range("O17").Value = 27.15485065307029

$ws.Range("B18").Value = 12.68314098929423
$ws.Range("C18").Value = 9.254636500206587
$ws.Range("E18").Value = 15.44727070812911
$ws.Range("F18").Value = 40.46764851204494
$ws.Range("G18").Value = 3.67687663825974
$ws.Range("I18").Value = 27.39843266498313
$ws.Range("J18").Value = 8.731425400091968
$ws.Range("K18").Value = 9.361951838325854
$ws.Range("L18").Value = 12.38895394540178
$ws.Range("M18").Value = 15.58442835713548
$ws.Range("O18").Value = 27.17555652551581

$ws.Range("B19").Value = 12.6628837007153
$ws.Range("C19").Value = 9.253841914794341
$ws.Range("E19").Value = 15.44772466525563
$ws.Range("F19").Value = 40.46994842411682
$ws.Range("G19").Value = 3.677007631180553
$ws.Range("I19").Value = 27.40615524550617
$ws.Range("J19").Value = 8.730609123450243
$ws.Range("K19").Value = 9.347534076298768
$ws.Range("L19").Value = 12.38792725954094
$ws.Range("M19").Value = 15.57965351517548
$ws.Range("O19").Value = 27.18264853168048

$ws.Range("B20").Value = 12.75388671879371
$ws.Range("C20").Value = 9.257423027150233
$ws.Range("E20").Value = 15.44580886862948
$ws.Range("F20").Value = 40.46015864569296
$ws.Range("G20").Value = 3.676421756249866
$ws.Range("I20").Value = 27.37169767553157
$ws.Range("J20").Value = 8.734280183635619
$ws.Range("K20").Value = 9.412302488676819
$ws.Range("L20").Value = 12.39264559461839
$ws.Range("M20").Value = 15.60125114823544
$ws.Range("O20").Value = 27.15105711036008

$ws.Range("B21").Value = 13.05662464338292
$ws.Range("C21").Value = 9.269547381256499
$ws.Range("E21").Value = 15.44159054388891
$ws.Range("F21").Value = 40.43706362750071
$ws.Range("G21").Value = 3.674516270144706
$ws.Range("I21").Value = 27.261085770289
$ws.Range("J21").Value = 8.746573298041088
$ws.Range("K21").Value = 9.62774896036176
$ws.Range("L21").Value = 12.41020629428653
$ws.Range("M21").Value = 15.67570926106005
$ws.Range("O21").Value = 27.05057757906414

$ws.Range("B22").Value = 13.25191249149335
$ws.Range("C22").Value = 9.277537697009912
$ws.Range("E22").Value = 15.44048066947556
$ws.Range("F22").Value = 40.42926472490579
$ws.Range("G22").Value = 3.673317663011207
$ws.Range("I22").Value = 27.1926410183379
$ws.Range("J22").Value = 8.754575278237848
$ws.Range("K22").Value = 9.766721693787812
$ws.Range("L22").Value = 12.42294398097692
$ws.Range("M22").Value = 15.72573661379275
$ws.Range("O22").Value = 26.98913403696726

$ws.Range("B23").Value = 13.14793858050295
$ws.Range("C23").Value = 9.273267202005353
$ws.Range("E23").Value = 15.44092275263982
$ws.Range("F23").Value = 40.43276021009707
$ws.Range("G23").Value = 3.673953118232052
$ws.Range("I23").Value = 27.22881930195222
$ws.Range("J23").Value = 8.750307528867989
$ws.Range("K23").Value = 9.692730823759147
$ws.Range("L23").Value = 12.41603107371713
$ws.Range("M23").Value = 15.69891421943533
$ws.Range("O23").Value = 27.02154066370816

$ws.Range("B24").Value = 12.74889979325324
$ws.Range("C24").Value = 9.257226016103491
$ws.Range("E24").Value = 15.44590571321912
$ws.Range("F24").Value = 40.46065938226752
$ws.Range("G24").Value = 3.676453692290326
$ws.Range("I24").Value = 27.373570505525
$ws.Range("J24").Value = 8.7340787372496
$ws.Range("K24").Value = 9.408753292847136
$ws.Range("L24").Value = 12.39238003181439
$ws.Range("M24").Value = 15.60005787192413
$ws.Range("O24").Value = 27.15277066660309

$ws.Range("B25").Value = 12.30954953418382
$ws.Range("C25").Value = 9.240223611044316
$ws.Range("E25").Value = 15.45839245937907
$ws.Range("F25").Value = 40.52211025480106
$ws.Range("G25").Value = 3.679352917194811
$ws.Range("I25").Value = 27.54622273776796
$ws.Range("J25").Value = 8.716445115958006
$ws.Range("K25").Value = 9.096013769689423
$ws.Range("L25").Value = 12.37235766465144
$ws.Range("M25").Value = 15.4995797558228
$ws.Range("O25").Value = 27.31243557639755
